# CW3MdigitalHandbook.docx - Additions to the chapter on the McKenzie Basin
# Wetlands Study.
#
# Adds a new "Baseline 2010-18 C98" scenario row (row 5) to the CW3M
# Clackamas regression-testing workbook, mirroring the existing baseline
# rows, with the columns that changed materially for this run (AET,
# basin discharge, and the two mass-balance-discrepancy columns)
# highlighted in yellow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 5 values -------------------------------------------------
$ws.Cells.Item(5, 1).Value = "CW3M"                      # A5
$ws.Cells.Item(5, 2).Value = "Baseline 2010-18 C98"      # B5
$ws.Cells.Item(5, 3).Value = 2010                        # C5

$ws.Cells.Item(5, 4).Value = 1143.6110839999999           # D5
$ws.Cells.Item(5, 5).Value = 1921.3682859999999           # E5
$ws.Cells.Item(5, 6).Value = 4.0370929999999996            # F5
$ws.Cells.Item(5, 7).Value = 197.01855499999999            # G5
$ws.Cells.Item(5, 8).Value = 73.459366000000003             # H5
$ws.Cells.Item(5, 9).Value = 122.410736                      # I5
$ws.Cells.Item(5, 10).Value = 62.789425000000001              # J5
$ws.Cells.Item(5, 11).Value = 717.20983899999999               # K5
$ws.Cells.Item(5, 12).Value = 75.975761000000006                # L5
$ws.Cells.Item(5, 13).Value = 999.48944100000006                 # M5
$ws.Cells.Item(5, 14).Value = 1219.0017089999999                  # N5
$ws.Cells.Item(5, 15).Value = 6236.0610349999997                   # O5
$ws.Cells.Item(5, 16).Value = 162867.046875                         # P5
$ws.Cells.Item(5, 17).Value = -387.43893400000002                    # Q5
$ws.Cells.Item(5, 18).Value = -0.111915                                # R5
$ws.Cells.Item(5, 19).Value = 2010                                      # S5

# D5:N5 share the plain 2-decimal number format used for the rest of the
# table (same as rows 2-4); O5:P5 use the integer format.
$ws.Range("D5:N5").NumberFormat = "0.00"
$ws.Range("O5:P5").NumberFormat = "0"
$ws.Range("Q5").NumberFormat = "0.00"
$ws.Range("R5").NumberFormat = "0.000000"

# Re-apply the 2-decimal format to the highlighted cells (NumberFormat
# must be set before the Interior color so the new yellow-filled styles
# get created with the right number format baked in).
$ws.Range("I5").NumberFormat = "0.00"
$ws.Range("M5").NumberFormat = "0.00"
$ws.Range("Q5").NumberFormat = "0.00"
$ws.Range("R5").NumberFormat = "0.000000"

# --- Highlight the cells that changed materially for this scenario ----
$ws.Range("I5").Interior.Color = 65535
$ws.Range("M5").Interior.Color = 65535
$ws.Range("Q5").Interior.Color = 65535
$ws.Range("R5").Interior.Color = 65535

# --- Sheet selection matches where the new scenario's delta columns are -
$ws.Range("Q5:R5").Select()
